# Rewrites the three "military / police / weapons" checkbox list items at
# the end of the TPS military-service addendum. They were numbered
# (ListParagraph / numId 1) bullet items; the new version drops the Word
# list numbering in favor of an inline {{ output_checkbox(...) }} merge
# field + a tab, using a first-line indent instead of list numbering, and
# moves the _GoBack bookmark from the end of the first item to the end of
# the last (now third) item.

$d = $word.ActiveDocument

# Locate the first and last paragraphs of the block to be rewritten by
# searching for stable anchor text rather than relying on fixed paragraph
# indices. wdFindContinue = 1, wdParagraph = 4 (used with Range.Expand so
# each Range covers its whole containing paragraph, not just the hit).
$startRange = $d.Content.Duplicate
$null = $startRange.Find.Execute("I was a member of the military.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startRange.Expand(4)

$endRange = $d.Content.Duplicate
$null = $endRange.Find.Execute("I received weapons", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange.Expand(4)

$blockRange = $d.Range($startRange.Start, $endRange.End)
$newXml = '<w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>output_checkbox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:r><w:t>users[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>].</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>served_military</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) }}</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">I was a member of the military. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>output_checkbox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:r><w:t>users[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>].</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>served_police</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) }}</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>I was a member of a police force.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>output_checkbox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:r><w:t>users[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>].</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>received_weapons_training</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and not </w:t></w:r><w:r><w:t>users[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>].</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>used_weapon_against_another_person</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) }}</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">I received weapons training, but I never used weapons against another person. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$blockRange.InsertXML($newXml)

Write-Host "Rewrote military/police/weapons checkbox block. Paragraph count now: $($d.Paragraphs.Count)"
